# Update stack-trace line numbers in the document body text to reflect
# the new line numbers after the M2Doc version was added to the
# template custom properties (issue #295).

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

Replace-Text "M2DocEvaluator.java:540)" "M2DocEvaluator.java:543)"
Replace-Text "M2DocEvaluator.java:1038)" "M2DocEvaluator.java:1084)"
Replace-Text "M2DocEvaluator.java:1254)" "M2DocEvaluator.java:1300)"
Replace-Text "M2DocEvaluator.java:275)" "M2DocEvaluator.java:278)"
Replace-Text "M2DocEvaluator.java:264)" "M2DocEvaluator.java:267)"
Replace-Text "M2DocUtils.java:712)" "M2DocUtils.java:694)"
Replace-Text "AbstractTemplatesTestSuite.java:459)" "AbstractTemplatesTestSuite.java:475)"
Replace-Text "AbstractTemplatesTestSuite.java:369)" "AbstractTemplatesTestSuite.java:384)"
